$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($text -eq "Think") {
        $p.Range.Text = "1. Think"
        $p.Style = "Heading 4"
    }
    elseif ($text -eq "Read") {
        $p.Range.Text = "2. Read"
        $p.Style = "Heading 4"
    }
}
